$d = $word.ActiveDocument

# The document ends with a (mostly empty) ListParagraph that only contains
# the "_GoBack" bookmark. The edit splits it into:
#   1) a new, still-empty ListParagraph
#   2) the original paragraph (keeping the bookmark), now carrying new text:
#      "This function converts text" + <bookmark> + " files containing
#       mapped reads into hic format files" + ". Users "

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()

# After the split, the bookmarked paragraph is now the last paragraph again.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)

# Append the two trailing runs at the end of the paragraph (after the
# bookmark, since the bookmark currently occupies the whole - empty - range).
$target.Range.InsertAfter(" files containing mapped reads into hic format files")
$target.Range.InsertAfter(". Users ")

# Insert the leading run before everything else in the paragraph (i.e.
# before the bookmark), by collapsing to the paragraph start.
$headRange = $target.Range.Duplicate
$headRange.Collapse(1)
$headRange.InsertBefore("This function converts text")
